# Update TPM-derived values on the active worksheet to reflect the new
# TPM computation (per commit message: "update scripts wuth new tpm").
# Only numeric values change; no structural/formula changes are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, taken from the canonical OOXML diff.
$updates = @{
    "G2"  = 6.206015333333333
    "H2"  = 18.618046
    "I2"  = 0.0150172404156507
    "J2"  = 0.0150172404156507
    "M2"  = 10.34761366666667
    "N2"  = 31.042841
    "O2"  = 0.2299953477621856
    "P2"  = 0.2299953477621856
    "Q2"  = 64.21744907874289
    "R2"  = 577.957041708686
    "S2"  = 0.00345389543182593
    "T2"  = 0.003453895431825931

    "G3"  = 6.206015333333333
    "H3"  = 18.618046
    "I3"  = 0.0150172404156507
    "J3"  = 0.0150172404156507
    "O3"  = 0.6794731949692173
    "P3"  = 0.6794731949692174
    "Q3"  = 189.7170343785558
    "R3"  = 1707.453309407002
    "S3"  = 0.01020381232484304
    "T3"  = 0.01020381232484304

    "G4"  = 6.206015333333333
    "H4"  = 18.618046
    "I4"  = 0.0150172404156507
    "J4"  = 0.0150172404156507
    "M4"  = 4.073058666666666
    "N4"  = 12.219176
    "O4"  = 0.09053145726859702
    "P4"  = 0.09053145726859703
    "Q4"  = 25.27746453889955
    "R4"  = 227.497180850096
    "S4"  = 0.001359532658981729
    "T4"  = 0.001359532658981729

    "I5"  = 0.9317452840597572
    "J5"  = 0.9317452840597571
    "M5"  = 10.34761366666667
    "N5"  = 31.042841
    "O5"  = 0.2299953477621856
    "P5"  = 0.2299953477621856
    "Q5"  = 3984.374204404961
    "R5"  = 35859.36783964466
    "S5"  = 0.2142970806331002
    "T5"  = 0.2142970806331002

    "I6"  = 0.9317452840597572
    "J6"  = 0.9317452840597571
    "O6"  = 0.6794731949692173
    "P6"  = 0.6794731949692174
    "S6"  = 0.6330959450575842
    "T6"  = 0.6330959450575843

    "I7"  = 0.9317452840597572
    "J7"  = 0.9317452840597571
    "M7"  = 4.073058666666666
    "N7"  = 12.219176
    "O7"  = 0.09053145726859702
    "P7"  = 0.09053145726859703
    "Q7"  = 1568.341301412593
    "R7"  = 14115.07171271334
    "S7"  = 0.0843522583690727
    "T7"  = 0.0843522583690727

    "G8"  = 22.00088566666667
    "H8"  = 66.002657
    "I8"  = 0.05323747552459213
    "J8"  = 0.05323747552459213
    "M8"  = 10.34761366666667
    "N8"  = 31.042841
    "O8"  = 0.2299953477621856
    "P8"  = 0.2299953477621856
    "Q8"  = 227.6566652031707
    "R8"  = 2048.909986828537
    "S8"  = 0.01224437169725941
    "T8"  = 0.01224437169725941

    "G9"  = 22.00088566666667
    "H9"  = 66.002657
    "I9"  = 0.05323747552459213
    "J9"  = 0.05323747552459213
    "O9"  = 0.6794731949692173
    "P9"  = 0.6794731949692174
    "Q9"  = 672.5640460414065
    "R9"  = 6053.076414372659
    "S9"  = 0.03617343758679013
    "T9"  = 0.03617343758679013

    "G10" = 22.00088566666667
    "H10" = 66.002657
    "I10" = 0.05323747552459213
    "J10" = 0.05323747552459213
    "M10" = 4.073058666666666
    "N10" = 12.219176
    "O10" = 0.09053145726859702
    "P10" = 0.09053145726859703
    "Q10" = 89.61089803895909
    "R10" = 806.498082350632
    "S10" = 0.004819666240542592
    "T10" = 0.004819666240542592
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
